$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a new bullet "Simplificar RangeUtil e Parse para só usar
#    Value2" immediately before the "Criar comando para remover um
#    lançamento em definitivo." bullet, re-using that paragraph's
#    list formatting.
# ------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("Criar comando para remover um lançamento em definitivo.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$target.Collapse(1)            # wdCollapseStart
$target.InsertParagraphBefore()
$newParaPos = $target.Start    # start of the freshly inserted paragraph

# Work out that paragraph's 1-based index and re-fetch it straight
# from the document's Paragraphs collection (rather than through the
# Find range itself) - reading .Range.Text off a paragraph obtained
# via a Find-scoped range/Paragraphs collection is unreliable in this
# host, but Document.Paragraphs.Item(n) works fine.
$paragraphIndex = $d.Range(0, $newParaPos).Paragraphs.Count + 1
$newParagraph = $d.Paragraphs.Item($paragraphIndex)

# NOTE: a trailing sentinel character "X" is appended to the new text
# so that, once the bookmark below has to be planted at the position
# immediately following the real text, that position does not fall
# exactly on the paragraph's end boundary (doing so trips a
# Bookmarks.Add edge case in this host that snaps/expands the
# bookmark across the whole paragraph instead of leaving it
# collapsed). The sentinel is stripped right after the bookmark has
# been created.
$newText = "Simplificar RangeUtil e Parse para só usar Value2"
$newParagraph.Range.Text = $newText + "X"

# ------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the end of the "Acho que as
#    entidades..." paragraph to the end of the newly inserted
#    paragraph (i.e. right after its text, before the paragraph
#    mark) - matching the diff.
# ------------------------------------------------------------------
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$newParagraphStart = $d.Paragraphs.Item($paragraphIndex).Range.Start
$bookmarkPos = $newParagraphStart + $newText.Length
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Strip the sentinel character now that the bookmark is anchored.
$sentinelRange = $d.Range($bookmarkPos, $bookmarkPos + 1)
$sentinelRange.Delete()
